$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("WIC")
$ws.Cells.Item(9, 1).Value = 44110.361111111109
$ws.Cells.Item(9, 2).Value = 634.5
$ws.Cells.Item(9, 3).Value = 12.8
